$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 467, pushing the existing data
# (previously rows 467:530) down to rows 469:532.
$ws.Rows("467:468").Insert()

# New record #1 (new row 467): Apio / Americana (o) / Primera,
# 2023-08-16 (serial 45154), sold "$/caja 8 unidades".
$ws.Range("A467").Value = 10
$ws.Range("B467").Value = "Vega Modelo de Temuco"
$ws.Range("C467").Value = "La Araucanía"
$ws.Range("D467").Value = 45154
$ws.Range("E467").Value = 9
$ws.Range("F467").Value = 100112017
$ws.Range("G467").Value = "Apio"
$ws.Range("H467").Value = "Americana (o)"
$ws.Range("I467").Value = "Primera"
$ws.Range("J467").Value = 90
$ws.Range("K467").Value = 10000
$ws.Range("L467").Value = 10000
$ws.Range("M467").Value = 10000
$ws.Range("N467").Value = "$/caja 8 unidades"
$ws.Range("O467").Value = "Provincia del Elquí"
$ws.Range("P467").Value = 10000
$ws.Range("Q467").Value = 1
$ws.Range("R467").Value = "Hortaliza"

# New record #2 (new row 468): Apio / Americana (o) / Primera,
# 2023-08-16 (serial 45154), sold "$/docena de matas".
$ws.Range("A468").Value = 10
$ws.Range("B468").Value = "Vega Modelo de Temuco"
$ws.Range("C468").Value = "La Araucanía"
$ws.Range("D468").Value = 45154
$ws.Range("E468").Value = 9
$ws.Range("F468").Value = 100112017
$ws.Range("G468").Value = "Apio"
$ws.Range("H468").Value = "Americana (o)"
$ws.Range("I468").Value = "Primera"
$ws.Range("J468").Value = 60
$ws.Range("K468").Value = 8000
$ws.Range("L468").Value = 8000
$ws.Range("M468").Value = 8000
$ws.Range("N468").Value = "$/docena de matas"
$ws.Range("O468").Value = "Provincia del Elquí"
$ws.Range("P468").Value = 1333
$ws.Range("Q468").Value = 6
$ws.Range("R468").Value = "Hortaliza"
